# PRCS252 - Use Case Analysis.docx edit script
# Implements: "Updated customer to include login use case" plus the
# accompanying run-merge/run-split housekeeping visible in the target diff.

$d = $word.ActiveDocument

function Merge-ParagraphRuns($searchText) {
    # Finds the paragraph whose text equals $searchText (already merged logically)
    # and forces Word to coalesce its runs into a single run by doing a
    # self-replace over the whole paragraph extent.
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
    return $found
}

# ---------------------------------------------------------------------
# 1) Driver / "view current stock" paragraph: 3 runs -> 1 run
# ---------------------------------------------------------------------
Merge-ParagraphRuns("As a driver, I want to be able to view current stock, so that the number of coaches available to drive in the depot is known.")

# ---------------------------------------------------------------------
# 2) Driver / "accept a e-ticket ... validate their journey" paragraph:
#    split the " validate their journey" run into " validate" + " their journey"
# ---------------------------------------------------------------------
$p = $d.Content
$p.Find.ClearFormatting()
$p.Find.Execute("to validate their journey on the coach.") | Out-Null
$pFound = $p.Find.Found
if ($pFound) {
    $paraRange = $p.Paragraphs(1).Range
    $fullText = $paraRange.Text
    $needle = " validate"
    $idx = $fullText.IndexOf($needle)
    if ($idx -ge 0) {
        $start = $paraRange.Start + $idx
        $len = $needle.Length
        $subRng = $d.Range($start, $start + $len)
        $subRng.Font.Bold = 1
        $subRng2 = $d.Range($start, $start + $len)
        $subRng2.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 3) Customer / "create an account" paragraph: 3 runs -> 1 run
# ---------------------------------------------------------------------
Merge-ParagraphRuns("As a customer, I want to be able to create an account, so that I can book a place on a coach.")

# ---------------------------------------------------------------------
# 4) Insert new Customer bullet: "login to my account" after the
#    "create an account" bullet and before "select the day and time".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("As a customer, I want to be able to create an account, so that I can book a place on a coach.") | Out-Null
if ($rng.Find.Found) {
    $accountPara = $rng.Paragraphs(1)
    $accountPara.Range.InsertParagraphAfter()
    # Re-find the paragraph collection to locate the freshly inserted (empty) paragraph.
    $newPara = $accountPara.Next()
    $newPara.Range.Text = "As a customer, I want to be able to login to my account, so that I can access my account."
}

# ---------------------------------------------------------------------
# 5) Customer / "create bookings ... catch the coach" paragraph: 3 runs -> 1 run
# ---------------------------------------------------------------------
Merge-ParagraphRuns("As a customer, I want to be able to create bookings, so that I can catch the coach to my desired destination.")

# ---------------------------------------------------------------------
# 6) Customer / "know my starting station" paragraph: 5 runs -> 1 run
# ---------------------------------------------------------------------
Merge-ParagraphRuns("As a customer, I want to be able to know my starting station, so I can know where to board the coach.")

# ---------------------------------------------------------------------
# 7) Customer / "know my ending station" paragraph: 5 runs -> 1 run
# ---------------------------------------------------------------------
Merge-ParagraphRuns("As a customer, I want to be able to know my ending station, so I can know where to get off the coach.")

# ---------------------------------------------------------------------
# 8) Admin / "approve customer account requests ... log in to their accounts."
#    paragraph: merge the final two runs only.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$old = ", so that customers can log in to their accounts."
$rng.Find.Execute(", so that customers can ", $true, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null

# ---------------------------------------------------------------------
# 9) Admin / "As an admin, I want to retrieve the timetables of the routes"
#    paragraph: split "As an admin, I want to retrieve the timetables of the "
#    into "As" + " an admin, I want to retrieve the timetables of the "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("As an admin, I want to retrieve the timetables of the routes") | Out-Null
if ($rng.Find.Found) {
    $paraRange = $rng.Paragraphs(1).Range
    $fullText = $paraRange.Text
    $needle = "As"
    $idx = $fullText.IndexOf($needle)
    if ($idx -ge 0) {
        $start = $paraRange.Start + $idx
        $len = $needle.Length
        $subRng = $d.Range($start, $start + $len)
        $subRng.Font.Bold = 1
        $subRng2 = $d.Range($start, $start + $len)
        $subRng2.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 10) Manager section
# ---------------------------------------------------------------------

# 10a) Move the _GoBack bookmark from the drawing paragraph to the
#      "Manager:" heading paragraph (covering just the "Manager:" text,
#      not its trailing paragraph mark).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Manager:") | Out-Null
if ($rng.Find.Found) {
    $managerPara = $rng.Paragraphs(1)
    $mStart = $managerPara.Range.Start
    $mEnd = $managerPara.Range.End - 1
    $bmRange = $d.Range($mStart, $mEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# 10b) "As" + " a manager, I want to be able to assign a driver a coach..."
#      merge into a single run.
Merge-ParagraphRuns("As a manager, I want to be able to assign a driver a coach for the shift, so that the driver can perform his duties.")

# 10c) "As a manager, I want to be able to remove an employee's shift, so
#      that they no longer undertake that shift." split into two runs:
#      "...undertake" + " that shift."
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("As a manager, I want to be able to remove an employee") | Out-Null
if ($rng.Find.Found) {
    $paraRange = $rng.Paragraphs(1).Range
    $fullText = $paraRange.Text
    $needle = " that shift."
    $idx = $fullText.LastIndexOf($needle)
    if ($idx -ge 0) {
        $start = $paraRange.Start + $idx
        $len = $needle.Length
        $subRng = $d.Range($start, $start + $len)
        $subRng.Font.Bold = 1
        $subRng2 = $d.Range($start, $start + $len)
        $subRng2.Font.Bold = 0
    }
}

Write-Output "done"
